$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for numeric-looking values so Excel
# does not silently coerce them into Number cells (which would
# strip significant trailing/leading zeros).
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D16", "D18", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "68.009.10"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "3.877.04"
$ws.Range("E3").Value = "  -0.77%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "481.23"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").Value = "144.02"
$ws.Range("E6").Value = "  -2.48%  "

$ws.Range("D7").Value = "0.618"
$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "0.733"
$ws.Range("E9").Value = "  +1.36%  "

$ws.Range("D10").Value = "0.178"
$ws.Range("E10").Value = "  +7.52%  "

$ws.Range("D11").Value = "0.0000350"
$ws.Range("E11").Value = "  -0.49%  "

$ws.Range("D12").Value = "42.67"
$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("D13").Value = "10.51"
$ws.Range("E13").Value = "  +1.96%  "

$ws.Range("D14").Value = "4.505.29"
$ws.Range("E14").Value = "  -0.35%  "

$ws.Range("D15").Value = "3.875.74"
$ws.Range("E15").Value = "  -1.71%  "

$ws.Range("D16").Value = "14.23"
$ws.Range("E16").Value = "  -2.65%  "

$ws.Range("E17").Value = "  -0.64%  "

$ws.Range("D18").Value = "20.07"
$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("E19").Value = "  -0.42%  "

$ws.Range("D20").Value = "68.085.45"
$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").Value = "427.82"
$ws.Range("E21").Value = "  -1.16%  "

$ws.Range("D22").Value = "3.52"
$ws.Range("E22").Value = "  +4.06%  "

$ws.Range("D23").Value = "14.72"
$ws.Range("E23").Value = "  +2.01%  "

$ws.Range("D24").Value = "89.88"
$ws.Range("E24").Value = "  +3.02%  "

$ws.Range("D25").Value = "11.92"
$ws.Range("E25").Value = "  +10.56%  "

$ws.Range("D26").Value = "3.66"
$ws.Range("E26").Value = "  +3.39%  "

$ws.Range("D27").Value = "11.00"
$ws.Range("E27").Value = "  +6.85%  "

$ws.Range("D28").Value = "37.29"
$ws.Range("E28").Value = "  -2.30%  "

$ws.Range("D29").Value = "5.67"
$ws.Range("E29").Value = "  -3.54%  "

$ws.Range("D30").Value = "715.70"
$ws.Range("E30").Value = "  -0.81%  "

$ws.Range("D31").Value = "13.51"
$ws.Range("E31").Value = "  +1.98%  "

$ws.Range("D32").Value = "0.129"
$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("D33").Value = "2.90"
$ws.Range("E33").Value = "  +2.89%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "6.07"
$ws.Range("E34").Value = "  +12.09%  "

$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0870"
$ws.Range("E35").Value = "  -2.12%  "

$ws.Range("D36").Value = "40.85"
$ws.Range("E36").Value = "  -2.07%  "

$ws.Range("D37").Value = "60.60"
$ws.Range("E37").Value = "  +3.08%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "0.997"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "0.395"
$ws.Range("E39").Value = "  +14.71%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.145"
$ws.Range("E40").Value = "  -4.03%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0497"
$ws.Range("E41").Value = "  +6.09%  "

$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "2.99"
$ws.Range("E42").Value = "  +4.68%  "

$ws.Range("D43").Value = "3.08"
$ws.Range("E43").Value = "  +3.69%  "

$ws.Range("D44").Value = "2.99"
$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.38"
$ws.Range("E45").Value = "  +4.19%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.142"
$ws.Range("E46").Value = "  +1.27%  "

$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("D48").Value = "3.37"
$ws.Range("E48").Value = "  -2.14%  "

$ws.Range("D49").Value = "2.11"
$ws.Range("E49").Value = "  -2.64%  "

$ws.Range("D50").Value = "144.70"
$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("D51").Value = "2.80"
$ws.Range("E51").Value = "  -1.74%  "

# Restore the default (General) style index on those cells so the
# only observable change is the cell text, matching the original
# (unstyled) cell formatting.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
